# Cost Estimation.xlsx - "Load cell" row update:
#   - cost 65 -> 200
#   - mark item as Purchased (checkmark in the Purchased column)
#   - the stray buydisplay.com link text in the link column is cleared
# Moving the cursor/selection to B26 (as last left by the author) and saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# "Load cell" row is row 14: B14=Load cell, C14=Cost, D14=Purchased?, F14=link
$ws.Range("C14").Value = 200

$ws.Range("D14").Value = "✅"
$ws.Range("D14").HorizontalAlignment = $xlCenter

$ws.Range("F14").ClearContents()

# Leave the selection where the author last left it before saving
$ws.Range("B26").Select()

$wb.Save()
